$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update the "Date" value (row 8, column B) ---
$ws.Range("B8").Value = "2021-12-22T21:26:07+01:00"

# --- Update the "Publisher" value (row 9, column B) ---
$ws.Range("B9").Value = "Forschungsgruppe Digital Health"

# --- Insert a new row 11 duplicating the "Contact" row (row 10), ---
# --- carrying over its formatting exactly via a range copy.      ---
$ws.Rows("11:11").Insert()
$ws.Range("A10:B10").Copy($ws.Range("A11:B11"))
$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "No display for ContactDetail"

# --- Insert a new row 12 for "Jurisdiction" / "Germany", again   ---
# --- copying formatting from the row above (now row 11).         ---
$ws.Rows("12:12").Insert()
$ws.Range("A11:B11").Copy($ws.Range("A12:B12"))
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = "Germany"
